# KIBON-2489 add new attributes to ferienbetreuung report
#
# Inserts three new columns (CC:CE) into the "Data" sheet, carrying the new
# "Delegationsmodell" attribute group (4_Sockelbeitrag,
# 4_Beiträge_nach_Anmeldungen, 4_davon_vorfinanzierte_Kantonsbeiträge), and
# splits the former "Kosten und Einnahmen" merged header so the trailing
# columns (Kantonsbeitrag / Kantonsbeitrag_anb / Gemeindebeteiligung) become
# their own "Resultate" group, while Kommentar stays un-grouped just like
# before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- 1. Make room: insert 3 new columns before the old CC (shifts
#        everything from CC onward, including cell content, merges,
#        column widths and the used-range dimension, to the right) ---
$ws.Range("CC:CE").Insert()

# Give the 3 new columns a sane custom width (same ballpark as their
# Resultate/Kantonsbeitrag neighbours).
$ws.Range("CC:CE").ColumnWidth = 20.7

# --- 2. Row 6: split the merged "Kosten und Einnahmen" banner ---
# Insert() widened the old BU6:CE6 merge to BU6:CH6 - undo that and build
# the three separate banners: "Kosten und Einnahmen" (BU:CB),
# "Delegationsmodell" (CC:CE, new) and "Resultate" (CF:CH).
$ws.Range("BU6:CH6").UnMerge()

$ws.Range("BU6:CB6").Merge()

$ws.Range("F6").Copy($ws.Range("CF6"))
$ws.Range("CF6").Value = "Resultate"
$ws.Range("G6").Copy($ws.Range("CG6"))
$ws.Range("AA6").Copy($ws.Range("CH6"))
$ws.Range("CF6:CH6").Merge()

$ws.Range("F6").Copy($ws.Range("CC6"))
$ws.Range("CC6").Value = "Delegationsmodell"
$ws.Range("G6").Copy($ws.Range("CD6"))
$ws.Range("AA6").Copy($ws.Range("CE6"))
$ws.Range("CC6:CE6").Merge()

# --- 3. Row 7: attribute-name headers for the 3 new columns ---
$ws.Range("CC7").Value = "4_Sockelbeitrag"
$ws.Range("CD7").Value = "4_Beiträge_nach_Anmeldungen"
$ws.Range("CE7").Value = "4_davon_vorfinanzierte_Kantonsbeiträge"

# --- 4. Row 8: placeholder tokens for the 3 new columns ---
$ws.Range("CC8").Value = "{sockelbeitrag}"
$ws.Range("CD8").Value = "{beitraegeNachAnmeldungen}"
$ws.Range("CE8").Value = "{vorfinanzierteKantonsbeitraege}"
